$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8 (shifts old row 8 and everything below down by one)
$ws.Rows("8:8").Insert()

# Fill in the new parameter row
$ws.Range("C8").Value = "道路设计标准流量参考值"
$ws.Range("D8").Value = "q_standard(v/h)"
$ws.Range("E8").Value = 10000

# Update selection to match the target workbook state
$ws.Range("E9").Select()
